{"js": "// Update the answer cells of the \"two-digit \u00f7 one-digit\" practice table.\n// The table has 20 rows x 5 columns, but only every 4th row (0, 4, 8, 12,\n// 16) actually holds an answer in each of its 5 cells; the rows in\n// between are blank spacer rows. We address each answer cell by its\n// (row, column) position so that duplicate old values (e.g. \"20\u00f77=2, 6\"\n// appears twice but maps to two different new values) are handled\n// correctly.\n\nconst table = context.document.body.tables.getFirst();\n\nconst updates = [\n  // [row, col, newText]\n  [0, 0, \"68\u00f79=7, 5\"],\n  [0, 1, \"55\u00f79=6, 1\"],\n  [0, 2, \"71\u00f77=10, 1\"],\n  [0, 3, \"71\u00f78=8, 7\"],\n  [0, 4, \"91\u00f76=15, 1\"],\n\n  [4, 0, \"13\u00f73=4, 1\"],\n  [4, 1, \"17\u00f73=5, 2\"],\n  [4, 2, \"41\u00f78=5, 1\"],\n  [4, 3, \"14\u00f79=1, 5\"],\n  [4, 4, \"52\u00f73=17, 1\"],\n\n  [8, 0, \"73\u00f78=9, 1\"],\n  [8, 1, \"68\u00f72=34, 0\"],\n  [8, 2, \"38\u00f79=4, 2\"],\n  [8, 3, \"18\u00f74=4, 2\"],\n  [8, 4, \"79\u00f75=15, 4\"],\n\n  [12, 0, \"65\u00f74=16, 1\"],\n  [12, 1, \"57\u00f79=6, 3\"],\n  [12, 2, \"61\u00f75=12, 1\"],\n  [12, 3, \"67\u00f76=11, 1\"],\n  [12, 4, \"62\u00f72=31, 0\"],\n\n  [16, 0, \"34\u00f79=3, 7\"],\n  [16, 1, \"73\u00f77=10, 3\"],\n  [16, 2, \"42\u00f78=5, 2\"],\n  [16, 3, \"71\u00f74=17, 3\"],\n  [16, 4, \"78\u00f76=13, 0\"],\n];\n\nfor (const [row, col, text] of updates) {\n  const cell = table.getCell(row, col);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the answer cells of the \"two-digit \u00f7 one-digit\" practice table.\n# The table has 20 rows x 5 columns, but only every 4th row (1, 5, 9, 13,\n# 17 in Word's 1-based Row numbering) actually holds an answer in each of\n# its 5 cells; the rows in between are blank spacer rows. We address each\n# answer cell by its 1-based (row, column) position so that duplicate old\n# values (e.g. \"20\u00f77=2, 6\" appears twice but maps to two different new\n# values) are handled correctly.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @(\n    @(1, 1, \"68\u00f79=7, 5\"),\n    @(1, 2, \"55\u00f79=6, 1\"),\n    @(1, 3, \"71\u00f77=10, 1\"),\n    @(1, 4, \"71\u00f78=8, 7\"),\n    @(1, 5, \"91\u00f76=15, 1\"),\n\n    @(5, 1, \"13\u00f73=4, 1\"),\n    @(5, 2, \"17\u00f73=5, 2\"),\n    @(5, 3, \"41\u00f78=5, 1\"),\n    @(5, 4, \"14\u00f79=1, 5\"),\n    @(5, 5, \"52\u00f73=17, 1\"),\n\n    @(9, 1, \"73\u00f78=9, 1\"),\n    @(9, 2, \"68\u00f72=34, 0\"),\n    @(9, 3, \"38\u00f79=4, 2\"),\n    @(9, 4, \"18\u00f74=4, 2\"),\n    @(9, 5, \"79\u00f75=15, 4\"),\n\n    @(13, 1, \"65\u00f74=16, 1\"),\n    @(13, 2, \"57\u00f79=6, 3\"),\n    @(13, 3, \"61\u00f75=12, 1\"),\n    @(13, 4, \"67\u00f76=11, 1\"),\n    @(13, 5, \"62\u00f72=31, 0\"),\n\n    @(17, 1, \"34\u00f79=3, 7\"),\n    @(17, 2, \"73\u00f77=10, 3\"),\n    @(17, 3, \"42\u00f78=5, 2\"),\n    @(17, 4, \"71\u00f74=17, 3\"),\n    @(17, 5, \"78\u00f76=13, 0\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $text = $u[2]\n    $cell = $t.Cell($row, $col)\n    $cell.Range.Text = $text\n}\n"}
